$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the mistyped password value ("secret_sauc" -> "secret_sauce") for all
# data rows so every account now uses the correct password.
$ws.Range("B3").Value = "secret_sauce"
$ws.Range("B4").Value = "secret_sauce"
$ws.Range("B5").Value = "secret_sauce"
$ws.Range("B6").Value = "secret_sauce"
$ws.Range("B7").Value = "secret_sauce"

# performance_glitch_user now correctly logs in (just slowly), so the
# expected result for that row becomes "Pass".
$ws.Range("C6").Value = "Pass"

# Update the remembered selection in the sheet view.
$ws.Range("F7").Select()
